# Title, label, URL, and organization columns working in Excel First Run.
#
# The "First run" sheet already has a D1 header "Organization"; this fills
# in the Organization column (D2:D5) with the same values as the Label
# column (B2:B5), and renames the placeholder "Sheet 1"/"Sheet 2" labels
# to "Org 1"/"Org 2" in both the Label and Organization columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First run")

# Rename the placeholder labels used for the two crawler sub-sheets.
$ws.Range("B3").Value = "Org 1"
$ws.Range("B4").Value = "Org 2"

# Populate the Organization column (D) to mirror the Label column (B).
$ws.Range("D2").Value = $ws.Range("B2").Value()
$ws.Range("D3").Value = $ws.Range("B3").Value()
$ws.Range("D4").Value = $ws.Range("B4").Value()
$ws.Range("D5").Value = $ws.Range("B5").Value()
